# Sample Project / Main.xlsx - "Rules" sheet edit.
# Row 11 (the "R40" rule row) had its B11 label cell changed from the
# text "R40" to the text "1". The new value must be stored as a text
# string (not a number), since "1" looks numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")

# Force the cell to text format first so the numeric-looking "1" is
# written out as a shared string (t="s") rather than coerced to a number.
$cell.NumberFormat = "@"
$cell.Value = "1"
